$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update "想去人数" (F column) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 8816
$ws1.Range("F3").Value = 8166
$ws1.Range("F9").Value = 146
$ws1.Range("F12").Value = 736
$ws1.Range("F14").Value = 4245
$ws1.Range("F17").Value = 17
$ws1.Range("F20").Value = 123

# Sheet "全部类型" (fourth sheet) - update "想去人数" (F column) values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8817
$ws4.Range("F3").Value = 8166
$ws4.Range("F9").Value = 146
$ws4.Range("F12").Value = 736
$ws4.Range("F14").Value = 4245
$ws4.Range("F17").Value = 17
$ws4.Range("F20").Value = 123
